$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New team order (row -> team abbreviation), reflecting the reshuffled
# shared-strings table from the commit ("fixed PER bug").
$teams = @{
    2  = "POR"
    3  = "NJN"
    4  = "CLE"
    5  = "DAL"
    6  = "ATL"
    7  = "OKC"
    8  = "CHA"
    9  = "WAS"
    10 = "MIL"
    11 = "LAC"
    12 = "SAS"
    13 = "DET"
    14 = "ORL"
    15 = "UTA"
    16 = "MEM"
    17 = "HOU"
    18 = "DEN"
    19 = "LAL"
    20 = "GSW"
    21 = "IND"
    22 = "CHI"
    23 = "PHI"
    24 = "BOS"
    25 = "TOR"
    26 = "MIA"
    27 = "SAC"
    28 = "PHO"
    29 = "NOH"
    30 = "NYK"
    31 = "MIN"
}

# New PER (per-minute) values replacing the old per-game point totals.
$values = @{
    2  = 13.32142857142857
    3  = 12.22142857142857
    4  = 11.87333333333333
    5  = 13.54615384615385
    6  = 11.88666666666666
    7  = 19.41666666666667
    8  = 10.40909090909091
    9  = 14.13846153846154
    10 = 13.02142857142857
    11 = 12.45384615384615
    12 = 14.91764705882353
    13 = 13.45833333333334
    14 = 13.8
    15 = 13.52
    16 = 13.23333333333334
    17 = 14.31538461538461
    18 = 13.22
    19 = 13.06153846153846
    20 = 13.76
    21 = 13.97857142857143
    22 = 14.175
    23 = 13.83076923076923
    24 = 14.33571428571429
    25 = 9.790909090909093
    26 = 12.025
    27 = 12.3625
    28 = 15.7
    29 = 12.82307692307692
    30 = 12.35
    31 = 11.81538461538462
}

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 2).Value = $teams[$row]
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
